$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two 2025-07-07 rows (old rows 2 and 3) - data for these sales
# was retired from the "vendas atipicas" extract. Remaining rows shift up.
$ws.Rows("2:3").Delete()

# The row that was previously row 5 (id_venda 379106, MATHEUS SILVEIRA /
# CAPA IPHONE 11) is now row 3; its estoque_atualizado value was corrected
# from 28 to 0.
$ws.Cells.Item(3, 7).Value = 0
